$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New GSW vs MIN playoff GM5 box-score rows (2025-05-14) appended as rows 24 & 25.
$rows = @(
    @{
        Row = 24
        A = 22; B = "GSW"; C = "MIN"; D = "away"; E = "2025-05-14"; F = "240:00";
        G = 39; H = 90; I = 0.433; J = 11; K = 39; L = 0.282; M = 21; N = 30; O = 0.7;
        P = 18; Q = 22; R = 40; S = 21; T = 14; U = 1; V = 17; W = 21; X = 110; Y = -11;
        Z = 23; AA = 24; AB = 25; AC = 38; AD = "L"
    },
    @{
        Row = 25
        A = 23; B = "MIN"; C = "GSW"; D = "home"; E = "2025-05-14"; F = "240:00";
        G = 49; H = 78; I = 0.628; J = 13; K = 31; L = 0.419; M = 10; N = 15; O = 0.667;
        P = 6; Q = 33; R = 39; S = 36; T = 11; U = 4; V = 20; W = 19; X = 121; Y = 11;
        Z = 30; AA = 32; AB = 31; AC = 28; AD = "W"
    }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Column A mirrors the bordered/bold "index" style used by the other data
    # rows (copy format only from an existing styled cell, e.g. A2).
    $ws.Range("A$rowNum").Value = $r.A
    $ws.Range("A2").Copy()
    $ws.Range("A$rowNum").PasteSpecial(-4122)

    $ws.Range("B$rowNum").Value = $r.B
    $ws.Range("C$rowNum").Value = $r.C
    $ws.Range("D$rowNum").Value = $r.D

    # Force the DATE column to stay plain text (matches the rest of the
    # sheet) instead of being auto-converted into a date serial number.
    $ws.Range("E$rowNum").NumberFormat = "@"
    $ws.Range("E$rowNum").Value = $r.E
    $ws.Range("E$rowNum").Style = "Normal"

    $ws.Range("F$rowNum").Value = $r.F
    $ws.Range("G$rowNum").Value = $r.G
    $ws.Range("H$rowNum").Value = $r.H
    $ws.Range("I$rowNum").Value = $r.I
    $ws.Range("J$rowNum").Value = $r.J
    $ws.Range("K$rowNum").Value = $r.K
    $ws.Range("L$rowNum").Value = $r.L
    $ws.Range("M$rowNum").Value = $r.M
    $ws.Range("N$rowNum").Value = $r.N
    $ws.Range("O$rowNum").Value = $r.O
    $ws.Range("P$rowNum").Value = $r.P
    $ws.Range("Q$rowNum").Value = $r.Q
    $ws.Range("R$rowNum").Value = $r.R
    $ws.Range("S$rowNum").Value = $r.S
    $ws.Range("T$rowNum").Value = $r.T
    $ws.Range("U$rowNum").Value = $r.U
    $ws.Range("V$rowNum").Value = $r.V
    $ws.Range("W$rowNum").Value = $r.W
    $ws.Range("X$rowNum").Value = $r.X
    $ws.Range("Y$rowNum").Value = $r.Y
    $ws.Range("Z$rowNum").Value = $r.Z
    $ws.Range("AA$rowNum").Value = $r.AA
    $ws.Range("AB$rowNum").Value = $r.AB
    $ws.Range("AC$rowNum").Value = $r.AC
    $ws.Range("AD$rowNum").Value = $r.AD
}
